$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '20.587.27'
$ws.Range('E2').Value = '  +1.90%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.473.71'
$ws.Range('E3').Value = '  +2.77%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.005'
$ws.Range('E4').Value = '  -0.35%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.9602'
$ws.Range('E5').Value = '  +5.39%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '277.67'
$ws.Range('E6').Value = '  +0.34%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.3628'
$ws.Range('E7').Value = '  -0.20%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3083'
$ws.Range('E8').Value = '  -0.36%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '39.75'
$ws.Range('E9').Value = '  +2.03%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.076'
$ws.Range('E10').Value = '  +6.02%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.06651'
$ws.Range('E11').Value = '  +2.27%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.000'
$ws.Range('E12').Value = '  -0.29%  '

$ws.Range('B13').Value = 'Solana'
$ws.Range('C13').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '18.29'
$ws.Range('E13').Value = '  +4.81%  '

$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.509'
$ws.Range('E14').Value = '  +3.09%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.9615'
$ws.Range('E15').Value = '  +1.97%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '6.167'
$ws.Range('E16').Value = '  +2.44%  '

$ws.Range('E17').Value = '  +1.42%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '1.474.05'
$ws.Range('E18').Value = '  +2.50%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.05945'
$ws.Range('E19').Value = '  +5.79%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '69.00'
$ws.Range('E20').Value = '  +2.22%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '5.507'
$ws.Range('E21').Value = '  +3.17%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '14.55'
$ws.Range('E22').Value = '  +1.96%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '11.19'
$ws.Range('E23').Value = '  +4.15%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.264'
$ws.Range('E24').Value = '  +1.13%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '20.583.07'
$ws.Range('E25').Value = '  +1.67%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '142.36'
$ws.Range('E26').Value = '  +4.14%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.135'
$ws.Range('E27').Value = '  -0.02%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '17.18'
$ws.Range('E28').Value = '  +2.09%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.635.40'
$ws.Range('E29').Value = '  +2.81%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '113.91'
$ws.Range('E30').Value = '  +4.26%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '3.910'
$ws.Range('E31').Value = '  +0.12%  '

$ws.Range('B32').Value = 'Stellar'
$ws.Range('C32').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.08023'
$ws.Range('E32').Value = '  +5.02%  '

$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.951'
$ws.Range('E33').Value = '  +3.83%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.8018'
$ws.Range('E34').Value = '  +1.12%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.506'
$ws.Range('E35').Value = '  +4.61%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.215'
$ws.Range('E36').Value = '  +7.55%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.05771'
$ws.Range('E37').Value = '  -2.59%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '4.744'
$ws.Range('E38').Value = '  +3.58%  '

$ws.Range('E39').Value = '  +3.88%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.9622'
$ws.Range('E40').Value = '  +4.42%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '10.42'
$ws.Range('E41').Value = '  +2.89%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.1881'
$ws.Range('E42').Value = '  +3.15%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '7.452'
$ws.Range('E43').Value = '  +6.39%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.5288'
$ws.Range('E44').Value = '  +1.83%  '

$ws.Range('B45').Value = 'PancakeSwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.523'
$ws.Range('E45').Value = '  +0.64%  '

$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '12.17'
$ws.Range('E46').Value = '  +1.88%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '118.81'
$ws.Range('E47').Value = '  +0.86%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.5219'
$ws.Range('E48').Value = '  +2.52%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.819'
$ws.Range('E49').Value = '  +4.43%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.06468'
$ws.Range('E50').Value = '  +2.65%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.9885'
$ws.Range('E51').Value = '  +0.20%  '
